$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67 (shifts existing rows 67-106 down to 68-107)
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new weekly record
$ws.Cells.Item(67, 1).Value = 10
$ws.Cells.Item(67, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(67, 3).Value = "La Araucanía"
$ws.Cells.Item(67, 4).Value = 44452
$ws.Cells.Item(67, 5).Value = 9
$ws.Cells.Item(67, 6).Value = 100114007
$ws.Cells.Item(67, 7).Value = "Jengibre"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 130
$ws.Cells.Item(67, 11).Value = 20000
$ws.Cells.Item(67, 12).Value = 25000
$ws.Cells.Item(67, 13).Value = 23077
$ws.Cells.Item(67, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(67, 15).Value = "Perú"
$ws.Cells.Item(67, 16).Value = 1775
$ws.Cells.Item(67, 17).Value = 13
$ws.Cells.Item(67, 18).Value = "Hortaliza"
